$d = $word.ActiveDocument

# Locate the paragraph that contains the text about the regression method,
# immediately following the inline picture added earlier in the document.
$target = $null
$blank = $null
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Для анализа данных был применен метод регрессии*") {
        $target = $para
        $blank = $d.Paragraphs.Item($i - 1)
        break
    }
}

if ($target -ne $null) {
    # Remove the paragraph with the commentary text, then the blank
    # paragraph that separated it from the picture.
    $target.Range.Delete()
    $blank.Range.Delete()
}
